$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555979631255"
$ws1.Range("B2").Value = "go_stims-16512555979321237.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255597946125.csv"
$ws1.Range("B4").Value = "go_stims-16512555979481306.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555979621267.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651255599151795"
$ws2.Range("B2").Value = "OB-16512555984497952.csv"
$ws2.Range("B3").Value = "ZB-match_3-16512555982287967.csv"
$ws2.Range("B4").Value = "OB-16512555983977957.csv"
$ws2.Range("B5").Value = "TB-16512555990747952.csv"
$ws2.Range("B6").Value = "OB-1651255598597794.csv"
$ws2.Range("B7").Value = "TB-16512555991347954.csv"
$ws2.Range("B8").Value = "ZB-match_4-1651255598281795.csv"
$ws2.Range("B9").Value = "TB-1651255598685798.csv"
$ws2.Range("B10").Value = "ZB-match_6-16512555981737947.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555991537955"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555991997964"
$ws4.Range("B2").Value = "MM_stims-1651255599167794.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555991558046.csv"
$ws4.Range("B4").Value = "MM_stims-16512555991827981.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555991687963.csv"
$ws4.Range("B6").Value = "MM_stims-16512555991987984.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555991838038.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1651255599278794"
$ws5.Range("B2").Value = "SAT_stims-1651255599205798.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651255599246796.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555992307959.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555992627969.csv"
